# Revised Budget.xlsx - add "ADDITIONS: made on 6/29/2017" line items to
# Sheet1 and add a new (empty) Sheet2, as per the commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Insert 5 new rows right above the current "Total" row (row 12),
#    pushing the Total row and the trailing blank rows down. The new
#    rows inherit row 11's D/E/F number formatting automatically.
# ---------------------------------------------------------------------
$ws1.Range("A12:F16").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2. Fill in the new "ADDITIONS" section header (row 11) and the five
#    new line items (rows 12-16).
# ---------------------------------------------------------------------
$ws1.Range("A11").Value = "ADDITIONS: made on 6/29/2017"
$ws1.Range("A11").Font.Bold = $true

$ws1.Range("A12").Value = "Ethernet Cables"
$ws1.Range("B12").Value = "Buhbo 6 inch CAT 7 Ethernet"
$ws1.Range("C12").Value = 1
$ws1.Range("D12").Value = 9.95
$ws1.Range("F12").Value = "https://www.amazon.com/Buhbo-Ethernet-Shielded-Network-Snagless/dp/B06XY8X7LY/ref=sr_1_3?ie=UTF8&qid=1498749192&sr=8-3&keywords=cat+7+6+inch+ethernet+cables"

$ws1.Range("A13").Value = "PoE Cable"
$ws1.Range("B13").Value = "BeElion PoE Injector and Splitter Kit"
$ws1.Range("C13").Value = 2
$ws1.Range("D13").Value = 8.09
$ws1.Range("F13").Value = "https://www.amazon.com/BeElion-Passive-Injector-Splitter-Connector/dp/B01HMNJHII/ref=sr_1_5?ie=UTF8&qid=1498743443&sr=8-5&keywords=PoE+Splitter"

$ws1.Range("A14").Value = "SMA Male to SMA Male"
$ws1.Range("B14").Value = "SMA Male to SMA Male Plug in series RF Coaxial Adapter Connector"
$ws1.Range("C14").Value = 3
$ws1.Range("D14").Value = 2.48
$ws1.Range("F14").Value = "https://www.amazon.com/Male-Plug-Coaxial-Adapter-Connector/dp/B007POCITA"

$ws1.Range("A15").Value = "N Female to SMA Female"
$ws1.Range("B15").Value = "N Female to SMA Female Adapter"
$ws1.Range("C15").Value = 3
$ws1.Range("D15").Value = 7.24
$ws1.Range("F15").Value = "https://www.showmecables.com/n-female-to-sma-female-adapter?gclid=CL6z_4La49QCFcqLswodVbYLRA"

$ws1.Range("A16").Value = "SMA Female to RP-SMA"
$ws1.Range("B16").Value = "DHT Electronics 2PCS RF coaxial coax adapter SMA female to RP-SMA male"
$ws1.Range("C16").Value = 1
$ws1.Range("D16").Value = 5.5
$ws1.Range("F16").Value = "https://www.amazon.com/DHT-Electronics-coaxial-adapter-female/dp/B00CVQ3XLY/ref=pd_bxgy_147_img_3?_encoding=UTF8&pd_rd_i=B00CVQ3XLY&pd_rd_r=E05PRYBB5DV00YT99V44&pd_rd_w=1QZRz&pd_rd_wg=5xQa4&psc=1&refRID=E05PRYBB5DV00YT99V44"

# Subtotal formulas for the new rows (shared formula C*D, same pattern as
# the rest of the sheet).
$ws1.Range("E12:E16").Formula = "=C12*D12"

# ---------------------------------------------------------------------
# 3. Fix up the grand Total formula (now on row 17) to include the new
#    rows.
# ---------------------------------------------------------------------
$ws1.Range("E17").Formula = "=SUM(E2:E16)"

# ---------------------------------------------------------------------
# 4. Restore the active-cell selection shown in the saved file.
# ---------------------------------------------------------------------
$ws1.Range("E12").Select()

# ---------------------------------------------------------------------
# 5. Add the new, empty "Sheet2" tab after "Sheet1".
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("G54").Select()

# Leave Sheet1 as the active/selected sheet, matching tabSelected="1".
$ws1.Select()
$ws1.Range("E12").Select()
